# pandas_column_formats.xlsx: replace the roster table with the "Creating
# reports" version — new headers, new per-student data, drop the trailing
# "Number of Passengers" column, and restyle the header row (bold+underline,
# no fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:F1), left to right -------------------------------------
$ws.Range("A1").Value = "Last Name"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Assigned School"
$ws.Range("D1").Value = "Host Teacher"
$ws.Range("E1").Value = "Practicum Course"
$ws.Range("F1").Value = "Practicum Day/Time"

# --- Data body, written column by column (matches the source export order) -

# Column A: Last Name
$ws.Range("A2").Value = "Links"
$ws.Range("A3").Value = "Nosar"
$ws.Range("A4").Value = "Scott"

# Column B: First Name
$ws.Range("B2").Value = "Kevin"
$ws.Range("B3").Value = "Bryan"
$ws.Range("B4").Value = "Spencer"

# Column C: Assigned School
$ws.Range("C2").Value = "Hugh Mercer Elementary School"
$ws.Range("C3").Value = "Hugh Mercer Elementary School"
$ws.Range("C4").Value = "A. G. Wright Middle School"

# Column D: Host Teacher
$ws.Range("D2").Value = "Nosar, Cathy"
$ws.Range("D3").Value = "Nosar, Cathy"
$ws.Range("D4").Value = "Coleman, Henry"

# Column E: Practicum Course
$ws.Range("E2").Value = "Art"
$ws.Range("E3").Value = "Computer"
$ws.Range("E4").Value = "Band III"

# Column F: Practicum Day/Time
$ws.Range("F2").Value = "M/T: 7:30:52 AM - 3:30:52 AM"
$ws.Range("F3").Value = "M/T/W/Th/F: 7:30:03 AM - 3:30:03 AM"
$ws.Range("F4").Value = "M/T: 7:30:06 AM - 3:30:06 AM"

# --- Drop the now-unused trailing column (old "Number of Passengers...") ---
$ws.Columns.Item(7).Delete()

# --- Header formatting: add underline to the bold header font, clear fill --
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Font.Underline = 2
$headerRange.Interior.Pattern = -4142
